# edit.ps1 -- "Added First Two Slides"
#
# Applies:
#   1. slide2: title "Header 1" -> "Soccer Data Analysis";
#      content "## [1] 3" (Courier) -> author list (default font).
#   2. Two new slides appended (Title and Content layout):
#        slide3: "Research Questions and Potential Methods"
#        slide4: "Dataset and Variable Details"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the existing second slide (Header 1 / ## [1] 3 placeholder).
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$titleShape = $s2.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Soccer Data Analysis"

$bodyShape = $s2.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange
$bodyRange.Text = "Aaron Graff, Juan Gonzalez, Andrew Henderson, Cody Farris"
$bodyRange.Font.Name = ""

# ---------------------------------------------------------------------
# 2. Append two new "Title and Content" slides after the current two.
# ---------------------------------------------------------------------
$titleContentLayout = $s2.CustomLayout

# --- Slide 3: Research Questions and Potential Methods ---
$s3 = $p.Slides.AddSlide(3, $titleContentLayout)

$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Research Questions and Potential Methods"

$s3Body = $s3.Shapes.Item(2).TextFrame.TextRange
$s3Body.Text = "Can we identify lower league players, as well as specific statistics that are indicators for potential success in La Liga?`rUsing modeling, can we distinguish which statistics in particular are most heavily correlated (most important) to player success in La Liga?`rWhat are the most valuable Key Performance Indicators (KPIs) in La Liga, as well as other lower-level leagues in Spain? Do lower-level leagues require more defense, for example, than La Liga?"
$s3Body.Paragraphs(1).IndentLevel = 1
$s3Body.Paragraphs(2).IndentLevel = 1
$s3Body.Paragraphs(3).IndentLevel = 1

# --- Slide 4: Dataset and Variable Details ---
$s4 = $p.Slides.AddSlide(4, $titleContentLayout)

$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Dataset and Variable Details"

$s4Body = $s4.Shapes.Item(2).TextFrame.TextRange
$s4Body.Text = "We have two main datasets that we will be utilizing:`rStathead Soccer Data: Data consisting of field players statistics (non-goalies) since the 2017-2018 season (when advanced statistics started to be calculated). These statistics include:"
$s4Body.Paragraphs(1).IndentLevel = 1
$s4Body.Paragraphs(2).IndentLevel = 2

# A second, duplicate content placeholder (mirrors the source deck, which
# carries the same two paragraphs twice on this slide).
$s4Body2Shape = $s4.Shapes.Item(2).Duplicate()
$s4Body2 = $s4Body2Shape.Item(1).TextFrame.TextRange
Write-Output $s4.Shapes.Count
